$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: wrap a bare <w:p>...</w:p> fragment into the flat-OPC "WordOpenXML"
# envelope that Range.InsertXML() / Range.InsertXML expects, then inject it
# so that it replaces the content of $range in place.
# ---------------------------------------------------------------------------
function Set-RangeParagraphXml($range, [string]$paragraphXml) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $paragraphXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

function Find-ParagraphContaining([string]$needle) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like ("*" + $needle + "*")) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) Insert a new (struck-through) paragraph right before the "Mask2" bullet,
#    and drop the <w:lastRenderedPageBreak/> from that bullet's first run
#    (it now opens the newly inserted paragraph instead).
# ---------------------------------------------------------------------------
$maskPara = Find-ParagraphContaining "Mask2"
$maskPara.Range.InsertParagraphBefore() | Out-Null

# Re-locate the (still unique) "Mask2" paragraph; the freshly minted empty
# paragraph now sits directly above it.
$maskPara = Find-ParagraphContaining "Mask2"
$newPara = $maskPara.Previous()

$newParaXml = @'
<w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="cs"/><w:strike/><w:rtl/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">תבטל </w:t></w:r><w:r><w:rPr><w:strike/></w:rPr><w:t>Motion capture</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:strike/><w:rtl/></w:rPr><w:t xml:space="preserve"> עבור </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:strike/></w:rPr><w:t>PAS</w:t></w:r></w:p>
'@
Set-RangeParagraphXml $newPara.Range $newParaXml

$maskParaXml = @'
<w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Mask2</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t xml:space="preserve"> ו-</w:t></w:r><w:r><w:t>mask3</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t xml:space="preserve"> מוצגות לזמן ארוך מדי בטריילים הראשונים.</w:t></w:r></w:p>
'@
Set-RangeParagraphXml $maskPara.Range $maskParaXml

# ---------------------------------------------------------------------------
# 2) Merge the ", " and "תתקן פונט חתוך" runs into a single run.
# ---------------------------------------------------------------------------
$fontPara = Find-ParagraphContaining "תסדר פונט כתב יד"
$fontParaXml = @'
<w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t>תסדר פונט כתב יד</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t>, תתקן פונט חתוך</w:t></w:r></w:p>
'@
Set-RangeParagraphXml $fontPara.Range $fontParaXml

# ---------------------------------------------------------------------------
# 3) Strike-through the "hand" bullet and append the "craig" sentence.
# ---------------------------------------------------------------------------
$handPara = Find-ParagraphContaining "האם משתמשים ביד אחת או 2?"
$handParaXml = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="cs"/><w:strike/><w:rtl/></w:rPr><w:t>האם משתמשים ביד אחת או 2?</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:strike/><w:rtl/></w:rPr><w:t xml:space="preserve"> נראה שביד 1 לפי מאמרים של </w:t></w:r><w:r><w:rPr><w:strike/></w:rPr><w:t>craig</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:strike/><w:rtl/></w:rPr><w:t>.</w:t></w:r></w:p>
'@
Set-RangeParagraphXml $handPara.Range $handParaXml

Write-Output "done"
